$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = '68.941.37'
$ws.Cells.Item(2, 5).Value = '  +2.38%  '
$ws.Cells.Item(3, 4).Value = '3.734.04'
$ws.Cells.Item(3, 5).Value = '  -1.12%  '
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.18%  '
$ws.Cells.Item(5, 4).Value = '601.76'
$ws.Cells.Item(5, 5).Value = '  +1.65%  '
$ws.Cells.Item(6, 4).Value = '168.81'
$ws.Cells.Item(6, 5).Value = '  -1.75%  '
$ws.Cells.Item(7, 4).Value = '3.732.75'
$ws.Cells.Item(7, 5).Value = '  -1.24%  '
$ws.Cells.Item(8, 4).Value = '1.00'
$ws.Cells.Item(8, 5).Value = '  +0.01%  '
$ws.Cells.Item(9, 4).Value = '0.532'
$ws.Cells.Item(9, 5).Value = '  +2.44%  '
$ws.Cells.Item(10, 4).Value = '0.164'
$ws.Cells.Item(10, 5).Value = '  +3.85%  '
$ws.Cells.Item(11, 4).Value = '6.27'
$ws.Cells.Item(11, 5).Value = '  -0.27%  '
$ws.Cells.Item(12, 4).Value = '0.461'
$ws.Cells.Item(12, 5).Value = '  +0.18%  '
$ws.Cells.Item(13, 4).Value = '38.18'
$ws.Cells.Item(13, 5).Value = '  +1.20%  '
$ws.Cells.Item(14, 4).Value = '0.0000244'
$ws.Cells.Item(14, 5).Value = '  +0.24%  '
$ws.Cells.Item(15, 4).Value = '4.357.96'
$ws.Cells.Item(15, 5).Value = '  -0.90%  '
$ws.Cells.Item(16, 4).Value = '3.734.62'
$ws.Cells.Item(16, 5).Value = '  -0.80%  '
$ws.Cells.Item(17, 4).Value = '68.923.96'
$ws.Cells.Item(17, 5).Value = '  +2.27%  '
$ws.Cells.Item(18, 4).Value = '7.25'
$ws.Cells.Item(18, 5).Value = '  +2.08%  '
$ws.Cells.Item(19, 4).Value = '0.115'
$ws.Cells.Item(19, 5).Value = '  +0.25%  '
$ws.Cells.Item(20, 4).Value = '17.21'
$ws.Cells.Item(20, 5).Value = '  +6.71%  '
$ws.Cells.Item(21, 4).Value = '496.79'
$ws.Cells.Item(21, 5).Value = '  +2.04%  '
$ws.Cells.Item(22, 4).Value = '9.52'
$ws.Cells.Item(22, 5).Value = '  +4.08%  '
$ws.Cells.Item(23, 4).Value = '0.723'
$ws.Cells.Item(23, 5).Value = '  +0.22%  '
$ws.Cells.Item(24, 4).Value = '84.79'
$ws.Cells.Item(24, 5).Value = '  +0.90%  '
$ws.Cells.Item(25, 4).Value = '0.0000142'
$ws.Cells.Item(25, 5).Value = '  +3.06%  '
$ws.Cells.Item(26, 4).Value = '2.31'
$ws.Cells.Item(26, 5).Value = '  -1.83%  '
$ws.Cells.Item(27, 4).Value = '12.27'
$ws.Cells.Item(27, 5).Value = '  +0.42%  '
$ws.Cells.Item(28, 4).Value = '10.12'
$ws.Cells.Item(28, 5).Value = '  -0.50%  '
$ws.Cells.Item(29, 4).Value = '1.00'
$ws.Cells.Item(29, 5).Value = '  -0.15%  '
$ws.Cells.Item(30, 4).Value = '2.93'
$ws.Cells.Item(30, 5).Value = '  +0.69%  '
$ws.Cells.Item(31, 4).Value = '2.43'
$ws.Cells.Item(31, 5).Value = '  +1.14%  '
$ws.Cells.Item(32, 4).Value = '7.96'
$ws.Cells.Item(32, 5).Value = '  +3.27%  '
$ws.Cells.Item(33, 4).Value = '31.65'
$ws.Cells.Item(33, 5).Value = '  -2.03%  '
$ws.Cells.Item(34, 4).Value = '3.875.03'
$ws.Cells.Item(34, 5).Value = '  -0.76%  '
$ws.Cells.Item(35, 4).Value = '0.109'
$ws.Cells.Item(35, 5).Value = '  +1.30%  '
$ws.Cells.Item(36, 4).Value = '3.669.80'
$ws.Cells.Item(36, 5).Value = '  -1.08%  '
$ws.Cells.Item(37, 4).Value = '1.00'
$ws.Cells.Item(37, 5).Value = '  +0.30%  '
$ws.Cells.Item(38, 4).Value = '1.01'
$ws.Cells.Item(38, 5).Value = '  +0.99%  '
$ws.Cells.Item(39, 4).Value = '5.79'
$ws.Cells.Item(39, 5).Value = '  +1.18%  '
$ws.Cells.Item(40, 4).Value = '0.133'
$ws.Cells.Item(40, 5).Value = '  -1.41%  '
$ws.Cells.Item(41, 4).Value = '0.323'
$ws.Cells.Item(41, 5).Value = '  +0.04%  '
$ws.Cells.Item(42, 4).Value = '436.12'
$ws.Cells.Item(42, 5).Value = '  -3.45%  '
$ws.Cells.Item(43, 4).Value = '49.02'
$ws.Cells.Item(43, 5).Value = '  +0.41%  '
$ws.Cells.Item(44, 4).Value = '1.99'
$ws.Cells.Item(44, 5).Value = '  +1.05%  '
$ws.Cells.Item(45, 4).Value = '2.87'
$ws.Cells.Item(45, 5).Value = '  +0.68%  '
$ws.Cells.Item(46, 4).Value = '8.38'
$ws.Cells.Item(46, 5).Value = '  +1.73%  '
$ws.Cells.Item(47, 4).Value = '1.00'
$ws.Cells.Item(47, 5).Value = '  +0.01%  '
$ws.Cells.Item(48, 4).Value = '40.41'
$ws.Cells.Item(48, 5).Value = '  -1.82%  '
$ws.Cells.Item(49, 4).Value = '143.41'
$ws.Cells.Item(49, 5).Value = '  +2.96%  '
$ws.Cells.Item(50, 4).Value = '0.0352'
$ws.Cells.Item(50, 5).Value = '  +1.23%  '
$ws.Cells.Item(51, 4).Value = '2.751.62'
$ws.Cells.Item(51, 5).Value = '  -2.07%  '
